$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting
# the existing N/O/P columns (Late / Outstanding-heading / Outstanding)
# one column to the right, and copy the column width from the column
# immediately to the left (M), matching Excel's default "Insert" behavior.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab with K18 selected,
# and move selection away from the "Transactions" sheet.
$ws.Activate() | Out-Null
$ws.Range("K18").Select() | Out-Null
